$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.94432133333333
$ws.Range("H2").Value = 62.832964
$ws.Range("I2").Value = 0.7396577289668299
$ws.Range("J2").Value = 0.7396577289668298
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.150777
$ws.Range("N2").Value = 0.452331
$ws.Range("O2").Value = 0.009673301965868179
$ws.Range("P2").Value = 0.009673301965868179
$ws.Range("Q2").Value = 3.157921937676
$ws.Range("R2").Value = 28.421297439084
$ws.Range("S2").Value = 0.007154932563684429
$ws.Range("T2").Value = 0.007154932563684427

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.94432133333333
$ws.Range("H3").Value = 62.832964
$ws.Range("I3").Value = 0.7396577289668299
$ws.Range("J3").Value = 0.7396577289668298
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.30706533333333
$ws.Range("N3").Value = 36.921196
$ws.Range("O3").Value = 0.7895763895222844
$ws.Range("P3").Value = 0.7895763895222843
$ws.Range("Q3").Value = 257.7631310116604
$ws.Range("R3").Value = 2319.868179104944
$ws.Range("S3").Value = 0.584016279119882
$ws.Range("T3").Value = 0.5840162791198817

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.94432133333333
$ws.Range("H4").Value = 62.832964
$ws.Range("I4").Value = 0.7396577289668299
$ws.Range("J4").Value = 0.7396577289668298
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.129079333333333
$ws.Range("N4").Value = 9.387238
$ws.Range("O4").Value = 0.2007503085118475
$ws.Range("P4").Value = 0.2007503085118475
$ws.Range("Q4").Value = 65.53644303482578
$ws.Range("R4").Value = 589.827987313432
$ws.Range("S4").Value = 0.1484865172832636
$ws.Range("T4").Value = 0.1484865172832635

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.327094666666667
$ws.Range("H5").Value = 6.981284
$ws.Range("I5").Value = 0.08218235047311259
$ws.Range("J5").Value = 0.08218235047311258
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.150777
$ws.Range("N5").Value = 0.452331
$ws.Range("O5").Value = 0.009673301965868179
$ws.Range("P5").Value = 0.009673301965868179
$ws.Range("Q5").Value = 0.350872352556
$ws.Range("R5").Value = 3.157851173004
$ws.Range("S5").Value = 0.0007949746923912277
$ws.Range("T5").Value = 0.0007949746923912276

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.327094666666667
$ws.Range("H6").Value = 6.981284
$ws.Range("I6").Value = 0.08218235047311259
$ws.Range("J6").Value = 0.08218235047311258
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 12.30706533333333
$ws.Range("N6").Value = 36.921196
$ws.Range("O6").Value = 0.7895763895222844
$ws.Range("P6").Value = 0.7895763895222843
$ws.Range("Q6").Value = 28.63970609951823
$ws.Range("R6").Value = 257.757354895664
$ws.Range("S6").Value = 0.06488924356901524
$ws.Range("T6").Value = 0.06488924356901522

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.327094666666667
$ws.Range("H7").Value = 6.981284
$ws.Range("I7").Value = 0.08218235047311259
$ws.Range("J7").Value = 0.08218235047311258
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.129079333333333
$ws.Range("N7").Value = 9.387238
$ws.Range("O7").Value = 0.2007503085118475
$ws.Range("P7").Value = 0.2007503085118475
$ws.Range("Q7").Value = 7.28166382817689
$ws.Range("R7").Value = 65.534974453592
$ws.Range("S7").Value = 0.01649813221170613
$ws.Range("T7").Value = 0.01649813221170612

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.044818
$ws.Range("H8").Value = 15.134454
$ws.Range("I8").Value = 0.1781599205600575
$ws.Range("J8").Value = 0.1781599205600575
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.150777
$ws.Range("N8").Value = 0.452331
$ws.Range("O8").Value = 0.009673301965868179
$ws.Range("P8").Value = 0.009673301965868179
$ws.Range("Q8").Value = 0.760642523586
$ws.Range("R8").Value = 6.845782712274
$ws.Range("S8").Value = 0.001723394709792523
$ws.Range("T8").Value = 0.001723394709792523

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.044818
$ws.Range("H9").Value = 15.134454
$ws.Range("I9").Value = 0.1781599205600575
$ws.Range("J9").Value = 0.1781599205600575
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 12.30706533333333
$ws.Range("N9").Value = 36.921196
$ws.Range("O9").Value = 0.7895763895222844
$ws.Range("P9").Value = 0.7895763895222843
$ws.Range("Q9").Value = 62.086904720776
$ws.Range("R9").Value = 558.7821424869841
$ws.Range("S9").Value = 0.1406708668333872
$ws.Range("T9").Value = 0.1406708668333872

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.044818
$ws.Range("H10").Value = 15.134454
$ws.Range("I10").Value = 0.1781599205600575
$ws.Range("J10").Value = 0.1781599205600575
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.129079333333333
$ws.Range("N10").Value = 9.387238
$ws.Range("O10").Value = 0.2007503085118475
$ws.Range("P10").Value = 0.2007503085118475
$ws.Range("Q10").Value = 15.785635744228
$ws.Range("R10").Value = 142.070721698052
$ws.Range("S10").Value = 0.03576565901687779
$ws.Range("T10").Value = 0.03576565901687779

$wb.Save()
